# Rename the "Larvae" lifestage value to "Larva" throughout the data
# (per commit message: changed "Larvae" to "Larva").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the extent of the data
$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

# Column C holds the "Lifestage" values; replace every "Larvae" with "Larva"
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Larvae") {
        $cell.Value = "Larva"
    }
}
